$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows of daily scores for 2025-02-28
$newRows = @(
    @("2025-02-28", "sleep", $true, $false),
    @("2025-02-28", "activity", $true, $false),
    @("2025-02-28", "weekly_activity", $true, $false)
)

$startRow = 83
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Leading apostrophe forces Excel to treat the date-like string as text
    # instead of auto-converting it to a date serial number; resetting the
    # style afterwards avoids leaving a stray quote-prefix number format.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = "'" + $rowData[0]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
}
